$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2: "Продавець: Rozetka." -> "Rozetka." (strip "Seller: " prefix)
$ws.Range("D2").Value = "Rozetka."

# H2: "Код:  395460480" -> "395460480" (strip "Code: " prefix; keep as text, not a number)
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "395460480"
$ws.Range("H2").Style = "Normal"

# I2: "Відгуки  137" -> "137" (strip "Reviews " prefix; keep as text, not a number)
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "137"
$ws.Range("I2").Style = "Normal"

# Remove row 3, which was an exact duplicate of row 2's product data
$ws.Rows("3").Delete()
